# Updated capital structure database
# Apply updated capital-structure figures to rows 2 and 3 (Thailand telecom equipment companies)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($row in 2,3) {
    $ws.Range("D$row").Value = -0.106
    $ws.Range("E$row").Value = 0.117
    $ws.Range("G$row").Value = -0.0328537170263789
    $ws.Range("H$row").Value = -0.0328537170263789
    $ws.Range("I$row").Value = -0.1039279534988656
    $ws.Range("J$row").Value = -0.09413047361076866
    $ws.Range("K$row").Value = 9.25
    $ws.Range("L$row").Value = 0.2218225419664268
    $ws.Range("M$row").Value = 3.3966
    $ws.Range("N$row").Value = 0.02779541734860884
    $ws.Range("O$row").Value = 0.3672
    $ws.Range("P$row").Value = 3.3966
    $ws.Range("Q$row").Value = 0.02779541734860884
    $ws.Range("R$row").Value = 0.3672
    $ws.Range("S$row").Value = 0
    $ws.Range("T$row").Value = 0
    $ws.Range("U$row").Value = 6.22
    $ws.Range("V$row").Value = 0.05090016366612111
    $ws.Range("W$row").Value = 0.2006507592190889
    $ws.Range("X$row").Value = 0.09387730119235498
    $ws.Range("Y$row").Value = 0.106773458026734
    $ws.Range("Z$row").Value = 0.4675151089519339
    $ws.Range("AA$row").Value = -0.04400741862583565
    $ws.Range("AB$row").Value = 0.07959387028541649
    $ws.Range("AC$row").Value = -0.1236012889112521
    $ws.Range("AD$row").Value = 26.5
    $ws.Range("AE$row").Value = 10.11897830451348
    $ws.Range("AF$row").Value = 36.61897830451348
    $ws.Range("AG$row").Value = 30.39897830451348
    $ws.Range("AH$row").Value = 0.2305705444994215
    $ws.Range("AI$row").Value = 0.416508233042485
    $ws.Range("AJ$row").Value = 0.1992082689036874
    $ws.Range("AK$row").Value = 0.3720851708966119
    $ws.Range("AL$row").Value = 1.13
    $ws.Range("AM$row").Value = 0.6719999999999999
    $ws.Range("AN$row").Value = 1324.999999999999
    $ws.Range("AO$row").Value = -4.026548672566372
    $ws.Range("AP$row").Value = 1519.948915225673
    $ws.Range("AQ$row").Value = -6.770833333333334
}
